$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks (C1 / G1 just get re-pointed at already-existing shared strings,
#     no visible text change: C1 stays "Allergens", G1 stays "LeaveEmpty") ---
$ws.Range("C1").Value = "Allergens"
$ws.Range("G1").Value = "LeaveEmpty"

# --- Row 2: Cobb Salad ---
$ws.Range("B2").Value = "Romaine / Egg / Guacamole / Bacon / Cheddar / Tomato / Cucumber / Ranch Dressing"
$ws.Range("C2").Value = "needed"
$ws.Range("D2").Value = "needed"
$ws.Range("E2").Value = "needed"
$ws.Range("F2").Value = "needed"

# --- Row 3: Kale Caesar ---
$ws.Range("B3").Value = "needed"
$ws.Range("C3").Value = "needed"
$ws.Range("D3").Value = "needed"
$ws.Range("E3").Value = "needed"
$ws.Range("F3").Value = "needed"

# --- Row 4: House Salad ---
$ws.Range("B4").Value = "Romaine / Carrots / Tomato / Cucumber / Balsamic Dressing"
$ws.Range("C4").Value = "needed"
$ws.Range("D4").Value = "needed"
$ws.Range("E4").Value = "needed"
$ws.Range("F4").Value = "needed"

# --- Formatting: center the new placeholder / content cells ---
$ws.Range("B3:C3").WrapText = $false
$ws.Range("B3:C3").HorizontalAlignment = -4108
$ws.Range("B3:C3").VerticalAlignment = -4108

$ws.Range("C4").WrapText = $false
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108

$ws.Range("D2:F2").HorizontalAlignment = -4108
$ws.Range("D2:F2").VerticalAlignment = -4108

$ws.Range("D3:F3").HorizontalAlignment = -4108
$ws.Range("D3:F3").VerticalAlignment = -4108

$ws.Range("D4:F4").HorizontalAlignment = -4108
$ws.Range("D4:F4").VerticalAlignment = -4108

$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108

$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").VerticalAlignment = -4108

# --- Selection as left by the author ---
$ws.Range("C2:C4").Select

Write-Host "edit applied"
